{"js": "// Update the worksheet date and the 25 multiplication answers (5 rows x 5\n// columns) with newly generated problems/answers, per the commit's\n// \"output generated at aa3dc9e\" refresh.\nconst replacements = [\n  [\"2023-11-09 Thursday\", \"2023-11-10 Friday\"],\n\n  [\"22\u00d777=1694\", \"72\u00d741=2952\"],\n  [\"19\u00d797=1843\", \"87\u00d714=1218\"],\n  [\"42\u00d726=1092\", \"32\u00d739=1248\"],\n  [\"66\u00d755=3630\", \"34\u00d727=918\"],\n  [\"99\u00d762=6138\", \"98\u00d783=8134\"],\n\n  [\"18\u00d784=1512\", \"89\u00d795=8455\"],\n  [\"40\u00d797=3880\", \"65\u00d786=5590\"],\n  [\"53\u00d760=3180\", \"88\u00d712=1056\"],\n  [\"41\u00d745=1845\", \"66\u00d753=3498\"],\n  [\"45\u00d717=765\", \"59\u00d733=1947\"],\n\n  [\"75\u00d734=2550\", \"29\u00d762=1798\"],\n  [\"28\u00d716=448\", \"46\u00d754=2484\"],\n  [\"61\u00d744=2684\", \"66\u00d759=3894\"],\n  [\"74\u00d765=4810\", \"14\u00d758=812\"],\n  [\"35\u00d798=3430\", \"64\u00d723=1472\"],\n\n  [\"72\u00d771=5112\", \"97\u00d712=1164\"],\n  [\"12\u00d742=504\", \"59\u00d758=3422\"],\n  [\"65\u00d741=2665\", \"78\u00d789=6942\"],\n  [\"30\u00d764=1920\", \"97\u00d718=1746\"],\n  [\"49\u00d725=1225\", \"61\u00d749=2989\"],\n\n  [\"69\u00d788=6072\", \"76\u00d781=6156\"],\n  [\"18\u00d734=612\", \"26\u00d717=442\"],\n  [\"67\u00d782=5494\", \"22\u00d759=1298\"],\n  [\"98\u00d760=5880\", \"48\u00d772=3456\"],\n  [\"95\u00d765=6175\", \"76\u00d723=1748\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 multiplication answers (5 rows x 5\n# columns) with newly generated problems/answers, per the commit's\n# \"output generated at aa3dc9e\" refresh.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-11-09 Thursday\", \"2023-11-10 Friday\"),\n\n    @(\"22\u00d777=1694\", \"72\u00d741=2952\"),\n    @(\"19\u00d797=1843\", \"87\u00d714=1218\"),\n    @(\"42\u00d726=1092\", \"32\u00d739=1248\"),\n    @(\"66\u00d755=3630\", \"34\u00d727=918\"),\n    @(\"99\u00d762=6138\", \"98\u00d783=8134\"),\n\n    @(\"18\u00d784=1512\", \"89\u00d795=8455\"),\n    @(\"40\u00d797=3880\", \"65\u00d786=5590\"),\n    @(\"53\u00d760=3180\", \"88\u00d712=1056\"),\n    @(\"41\u00d745=1845\", \"66\u00d753=3498\"),\n    @(\"45\u00d717=765\", \"59\u00d733=1947\"),\n\n    @(\"75\u00d734=2550\", \"29\u00d762=1798\"),\n    @(\"28\u00d716=448\", \"46\u00d754=2484\"),\n    @(\"61\u00d744=2684\", \"66\u00d759=3894\"),\n    @(\"74\u00d765=4810\", \"14\u00d758=812\"),\n    @(\"35\u00d798=3430\", \"64\u00d723=1472\"),\n\n    @(\"72\u00d771=5112\", \"97\u00d712=1164\"),\n    @(\"12\u00d742=504\", \"59\u00d758=3422\"),\n    @(\"65\u00d741=2665\", \"78\u00d789=6942\"),\n    @(\"30\u00d764=1920\", \"97\u00d718=1746\"),\n    @(\"49\u00d725=1225\", \"61\u00d749=2989\"),\n\n    @(\"69\u00d788=6072\", \"76\u00d781=6156\"),\n    @(\"18\u00d734=612\", \"26\u00d717=442\"),\n    @(\"67\u00d782=5494\", \"22\u00d759=1298\"),\n    @(\"98\u00d760=5880\", \"48\u00d772=3456\"),\n    @(\"95\u00d765=6175\", \"76\u00d723=1748\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
